# Afegida funció formula_text() amb millores
#
# Populates Hoja1 with the iris-variable codebook (camp/etiqueta/formula/
# recode/rank), centers the numeric "formula" / "recode" / "rank" columns,
# sizes the first two label columns, and leaves the selection where the
# author left it (D17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- cells, written in the same order the original author entered them --
$ws.Range("A2").Value = "Sepal.Length"
$ws.Range("A1").Value = "camp"
$ws.Range("A4").Value = "Petal.Length"
$ws.Range("A5").Value = "Petal.Width"
$ws.Range("A6").Value = "Species"
$ws.Range("A3").Value = "Sepal.Width"
$ws.Range("B1").Value = "etiqueta"
$ws.Range("C1").Value = "formula"
$ws.Range("B2").Value = "Lengh Sepal"
$ws.Range("B3").Value = "Widht of Sepal"
$ws.Range("B4").Value = "lengh of petal "
$ws.Range("B5").Value = "Widht of Petal"
$ws.Range("B6").Value = "Species varies "
$ws.Range("D1").Value = "recode"
$ws.Range("E1").Value = "rank"
$ws.Range("A7").Value = "country"
$ws.Range("B7").Value = "Catalonia"
$ws.Range("A8").Value = "regiligio"
$ws.Range("B8").Value = "fsffadf"

# --- numeric "formula" column --------------------------------------------
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 4

# --- formatting ---------------------------------------------------------
# Center the "formula" / "recode" / "rank" columns (only the populated cells)
$ws.Range("C1:C8").HorizontalAlignment = -4108  # xlCenter
$ws.Range("D1:E1").HorizontalAlignment = -4108  # xlCenter

# Column widths for the label columns (17.7109375 / 16.42578125 chars on-disk,
# closest the rendering grid allows is 16.8 / 15.65 typed in the width box)
$ws.Columns.Item(1).ColumnWidth = 16.8
$ws.Columns.Item(2).ColumnWidth = 15.65

# --- selection, matching the author's last cursor position --------------
$ws.Range("D17").Select()
